# Applies the Result_Sheet.xlsx restructuring:
#  - Remove the Q1/Q2/Q3 columns (old E,F,G) so Total Marks/Percentage/Status
#    (old H,I,J) shift left into E,F,G
#  - Update the "Marks Obtained" totals in column D for the graded rows
#  - Fill in the previously "Not Graded" row (row 14) with real marks
#  - Restore the title merge across A1:F1
#  - Restore dimension / column widths (handled automatically by the
#    column delete, which preserves exact widths from the old H/I/J cols)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old Q1 (E), Q2 (F), Q3 (G) columns; Total Marks/Percentage/Status
# (old H/I/J) shift left to become the new E/F/G, carrying their widths
# (13, 12, 9) and values with them automatically.
$ws.Columns("E:G").Delete()

# Fix up "Marks Obtained" (column D) for the graded student rows - these
# values changed as part of the restructuring.
$ws.Range("D9").Value = 11.5
$ws.Range("D10").Value = 14
$ws.Range("D11").Value = 14
$ws.Range("D12").Value = 13
$ws.Range("D13").Value = 14

# Row 14 previously had no Total/Percentage/Status (student was "Not
# Graded"); populate it with final grading results.
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 20
# Force "30.0%" to be stored as literal text (matching the other Percentage
# cells) instead of letting Excel auto-convert it to a percentage number;
# ClearFormats afterwards drops the temporary text format so the cell keeps
# the default (unstyled) look, same as its siblings.
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "30.0%"
$ws.Range("F14").ClearFormats()
$ws.Range("G14").Value = "Checked"

# The column delete shrank the title merge (it used to span A1:F1); put it
# back. (Merging re-stamps every covered cell, so clear the formatting
# Excel applies to the newly-covered, still-empty B1:F1 cells.)
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:F1").Merge()
$ws.Range("B1:F1").ClearFormats()
